# Add Test data for Spain Zettler Market
# - Duplicate the "Italy" sheet (same layout/styles/merges) placing the
#   copy right after it, rename the copy to "Spain".
# - Fill in the Spain-specific market name and part-number code.
# - Bump the row heights of the "used for / constants / user story" block
#   (rows 3-5) and widen the columns to fit the longer Spain text.
# - Select B11:B12 on the new Spain sheet (becomes the active sheet/tab),
#   and reset Italy's selection back to the full used range A1:D19.

$wb = $excel.ActiveWorkbook

$italy = $wb.Worksheets.Item("Italy")

# Duplicate Italy (carries over styles, merged cells, shared-string
# references, page setup, ...) and place the new sheet right after it.
$italy.Copy([System.Reflection.Missing]::Value, $italy)
$spain = $wb.Worksheets.Item($italy.Index + 1)
$spain.Name = "Spain"

# New market-specific cell content.
# B4 is written before B2 so the new shared strings land in the same
# order as the target workbook (NGC code first, then "Spain Market").
$spain.Range("B4").Value = "NGC-3103/T2019/T2066/T2046"
$spain.Range("B2").Value = "Spain Market"

# Taller rows for the wrapped-text instructions in column D.
$spain.Rows.Item(3).RowHeight = 28.8
$spain.Rows.Item(4).RowHeight = 28.8
$spain.Rows.Item(5).RowHeight = 28.8

# Wider columns to fit the Spain market values.
$spain.Columns.Item(1).ColumnWidth = 25.109375
$spain.Columns.Item(2).ColumnWidth = 27.109375
$spain.Columns.Item(3).ColumnWidth = 14.109375
$spain.Columns.Item(4).ColumnWidth = 18.109375

# Reset Italy's selection (it's no longer the active/selected tab).
$italy.Range("A1:D19").Select()

# Select B11:B12 on Spain last so it ends up the active sheet/tab.
$spain.Range("B11:B12").Select()
